$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2142.375
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 2319.8572
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 2319.8572
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -2457.8572
$ws.Range("H125").Value = 4354.6294
$ws.Range("I125").Value = 477.33334
$ws.Range("J125").Value = 4839.2915
$ws.Range("K125").Value = 4296.00006
$ws.Range("L125").Value = 43553.6235
$ws.Range("M125").Value = -1836.00006
$ws.Range("N125").Value = -48473.6235
$ws.Range("H129").Value = 1099.6389
$ws.Range("J129").Value = 1132.9565
$ws.Range("L129").Value = 3398.8695
$ws.Range("N129").Value = -13398.8695
$ws.Range("H132").Value = 2102.6875
$ws.Range("I132").Value = 2223.7932
$ws.Range("J132").Value = 932
$ws.Range("K132").Value = 6671.3796
$ws.Range("L132").Value = 2796
$ws.Range("M132").Value = -4141.3796
$ws.Range("N132").Value = -7856
$ws.Range("H137").Value = 2970
$ws.Range("I137").Value = 3069.7273
$ws.Range("K137").Value = 9209.1819
$ws.Range("M137").Value = -6659.1819
$ws.Range("H138").Value = 2560.8157
$ws.Range("I138").Value = 1015.88
$ws.Range("K138").Value = 3047.64
$ws.Range("M138").Value = 2092.36
$ws.Range("H141").Value = 1757.7446
$ws.Range("I141").Value = 1393.3096
$ws.Range("J141").Value = 4819
$ws.Range("K141").Value = 4179.9288
$ws.Range("L141").Value = 14457
$ws.Range("M141").Value = 1000.0712
$ws.Range("N141").Value = -24817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1229.409
$ws.Range("I2").Value = 1431.6428
$ws.Range("J2").Value = 875.5
$ws.Range("K2").Value = 1431.6428
$ws.Range("L2").Value = 875.5
$ws.Range("M2").Value = -1318.6428
$ws.Range("N2").Value = -1101.5
$ws.Range("H32").Value = 23509.91
$ws.Range("I32").Value = 25605.262
$ws.Range("J32").Value = 13871.3
$ws.Range("K32").Value = 25605.262
$ws.Range("L32").Value = 13871.3
$ws.Range("M32").Value = -25318.262
$ws.Range("N32").Value = -14445.3
$ws.Range("H97").Value = 1238.4615
$ws.Range("I97").Value = 1310
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1310
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -814
$ws.Range("N97").Value = -1992
$ws.Range("H110").Value = 1071.96
$ws.Range("I110").Value = 1068.9546
$ws.Range("J110").Value = 1094
$ws.Range("K110").Value = 1068.9546
$ws.Range("L110").Value = 1094
$ws.Range("M110").Value = 976.0454
$ws.Range("N110").Value = -5184
$ws.Range("H116").Value = 1229.409
$ws.Range("I116").Value = 1431.6428
$ws.Range("J116").Value = 875.5
$ws.Range("K116").Value = 1431.6428
$ws.Range("L116").Value = 875.5
$ws.Range("M116").Value = 862.3571999999999
$ws.Range("N116").Value = -5463.5
$ws.Range("H122").Value = 2414.5881
$ws.Range("I122").Value = 2651.4285
$ws.Range("J122").Value = 2248.8
$ws.Range("K122").Value = 7954.2855
$ws.Range("L122").Value = 6746.400000000001
$ws.Range("M122").Value = -5504.2855
$ws.Range("N122").Value = -11646.4
$ws.Range("H132").Value = 3933.5085
$ws.Range("I132").Value = 1800.55
$ws.Range("J132").Value = 8423.947
$ws.Range("K132").Value = 5401.65
$ws.Range("L132").Value = 25271.841
$ws.Range("M132").Value = -2871.65
$ws.Range("N132").Value = -30331.841

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1229.409
$ws.Range("I3").Value = 1431.6428
$ws.Range("J3").Value = 875.5
$ws.Range("K3").Value = 1431.6428
$ws.Range("L3").Value = 875.5
$ws.Range("M3").Value = -1317.6428
$ws.Range("N3").Value = -1103.5
$ws.Range("H13").Value = 69150
$ws.Range("J13").Value = 69150
$ws.Range("L13").Value = 69150
$ws.Range("N13").Value = -69486
$ws.Range("H80").Value = 197.33333
$ws.Range("I80").Value = 147.5
$ws.Range("J80").Value = 201.86363
$ws.Range("K80").Value = 147.5
$ws.Range("L80").Value = 201.86363
$ws.Range("M80").Value = 850.5
$ws.Range("N80").Value = -2197.86363
$ws.Range("H83").Value = 197.33333
$ws.Range("I83").Value = 147.5
$ws.Range("J83").Value = 201.86363
$ws.Range("K83").Value = 737.5
$ws.Range("L83").Value = 1009.31815
$ws.Range("M83").Value = 4254.5
$ws.Range("N83").Value = -10993.31815
$ws.Range("H94").Value = 1121.125
$ws.Range("I94").Value = 948.9091
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 948.9091
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -497.9091
$ws.Range("N94").Value = -2402
$ws.Range("H105").Value = 2088960.8
$ws.Range("I105").Value = 7820750
$ws.Range("J105").Value = 4673.8184
$ws.Range("K105").Value = 7820750
$ws.Range("L105").Value = 4673.8184
$ws.Range("M105").Value = -7819003
$ws.Range("N105").Value = -8167.8184
$ws.Range("H107").Value = 1287.091
$ws.Range("I107").Value = 1214.5
$ws.Range("J107").Value = 2013
$ws.Range("K107").Value = 1214.5
$ws.Range("L107").Value = 2013
$ws.Range("M107").Value = 705.5
$ws.Range("N107").Value = -5853

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1267.6666
$ws.Range("I16").Value = 781.7143
$ws.Range("J16").Value = 1467.7646
$ws.Range("K16").Value = 781.7143
$ws.Range("L16").Value = 1467.7646
$ws.Range("M16").Value = -494.7143
$ws.Range("N16").Value = -2041.7646
$ws.Range("H31").Value = 4949.3076
$ws.Range("I31").Value = 5626.64
$ws.Range("J31").Value = 3739.7856
$ws.Range("K31").Value = 5626.64
$ws.Range("L31").Value = 3739.7856
$ws.Range("M31").Value = -5331.64
$ws.Range("N31").Value = -4329.7856
$ws.Range("H34").Value = 4949.3076
$ws.Range("I34").Value = 5626.64
$ws.Range("J34").Value = 3739.7856
$ws.Range("K34").Value = 5626.64
$ws.Range("L34").Value = 3739.7856
$ws.Range("M34").Value = -5424.64
$ws.Range("N34").Value = -4143.7856
$ws.Range("H58").Value = 2068799.9
$ws.Range("I58").Value = 3498497.2
$ws.Range("J58").Value = 3681.6667
$ws.Range("K58").Value = 3498497.2
$ws.Range("L58").Value = 3681.6667
$ws.Range("M58").Value = -3498294.2
$ws.Range("N58").Value = -4087.6667
$ws.Range("H113").Value = 1267.6666
$ws.Range("I113").Value = 781.7143
$ws.Range("J113").Value = 1467.7646
$ws.Range("K113").Value = 781.7143
$ws.Range("L113").Value = 1467.7646
$ws.Range("M113").Value = 1388.2857
$ws.Range("N113").Value = -5807.7646
$ws.Range("H132").Value = 3995.157
$ws.Range("I132").Value = 4071.375
$ws.Range("J132").Value = 3718
$ws.Range("K132").Value = 12214.125
$ws.Range("L132").Value = 11154
$ws.Range("M132").Value = -9684.125
$ws.Range("N132").Value = -16214
$ws.Range("H136").Value = 2068799.9
$ws.Range("I136").Value = 3498497.2
$ws.Range("J136").Value = 3681.6667
$ws.Range("K136").Value = 10495491.6
$ws.Range("L136").Value = 11045.0001
$ws.Range("M136").Value = -10492941.6
$ws.Range("N136").Value = -16145.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3044.8572
$ws.Range("I124").Value = 865
$ws.Range("J124").Value = 3916.8
$ws.Range("K124").Value = 2595
$ws.Range("L124").Value = 11750.4
$ws.Range("M124").Value = 2315
$ws.Range("N124").Value = -21570.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7856.2573
$ws.Range("I132").Value = 5160.2256
$ws.Range("J132").Value = 28750.5
$ws.Range("K132").Value = 15480.6768
$ws.Range("L132").Value = 86251.5
$ws.Range("M132").Value = -12950.6768
$ws.Range("N132").Value = -91311.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1119.174
$ws.Range("I16").Value = 656.41174
$ws.Range("J16").Value = 2430.3333
$ws.Range("K16").Value = 656.41174
$ws.Range("L16").Value = 2430.3333
$ws.Range("M16").Value = -486.41174
$ws.Range("N16").Value = -2770.3333
$ws.Range("H46").Value = 1749.1818
$ws.Range("I46").Value = 1956.2
$ws.Range("J46").Value = 1576.6666
$ws.Range("K46").Value = 1956.2
$ws.Range("L46").Value = 1576.6666
$ws.Range("M46").Value = -1768.2
$ws.Range("N46").Value = -1952.6666
$ws.Range("H122").Value = 6035.1143
$ws.Range("I122").Value = 5747.891
$ws.Range("J122").Value = 6585.625
$ws.Range("K122").Value = 17243.673
$ws.Range("L122").Value = 19756.875
$ws.Range("M122").Value = -14793.673
$ws.Range("N122").Value = -24656.875
$ws.Range("H132").Value = 6869.619
$ws.Range("I132").Value = 8886.68
$ws.Range("J132").Value = 3903.353
$ws.Range("K132").Value = 26660.04
$ws.Range("L132").Value = 11710.059
$ws.Range("M132").Value = -24130.04
$ws.Range("N132").Value = -16770.059
$ws.Range("H136").Value = 3695.6826
$ws.Range("I136").Value = 2211
$ws.Range("J136").Value = 5436.3447
$ws.Range("K136").Value = 6633
$ws.Range("L136").Value = 16309.0341
$ws.Range("M136").Value = -4083
$ws.Range("N136").Value = -21409.0341

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2393.6365
$ws.Range("I132").Value = 1268.875
$ws.Range("J132").Value = 3452.2354
$ws.Range("K132").Value = 3806.625
$ws.Range("L132").Value = 10356.7062
$ws.Range("M132").Value = -1276.625
$ws.Range("N132").Value = -15416.7062
$ws.Range("H136").Value = 7057.7827
$ws.Range("I136").Value = 7154.5
$ws.Range("J136").Value = 6983.385
$ws.Range("K136").Value = 21463.5
$ws.Range("L136").Value = 20950.155
$ws.Range("M136").Value = -18913.5
$ws.Range("N136").Value = -26050.155

Write-Host "Done applying changes"